# Leave Card update — 4/12/2023 4:43 PM
#
# A new "SOLO P(1-0-0)" leave entry (1.25 days, dated 3/6/2023) is inserted
# into Table1 right after the existing row for period 2/21/2023
# ("SP(1-0-0)" / "BDAY 3/7/23", row 521). That existing row also gets its
# EARNED value filled in (1.25), which the calculated BALANCE columns
# (row 9) pick up automatically on recalculation. Inserting the row shifts
# every following table row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item(1)

# 1) Fill in the EARNED value for the existing "SP(1-0-0)" row (521).
#    Its [EARNED ] calculated column (G) recomputes to 1.25 automatically.
$ws.Cells.Item(521, 3).Value = 1.25

# 2) Insert a blank row at 522, pushing the old rows 522:705 down to 523:706.
$ws.Rows.Item(522).Insert()

# 3) Grow the table so it covers the newly inserted row.
$tbl.Resize($ws.Range("A8:K706"))

# 4) The row that now sits at the bottom of the resized table (706, the old
#    705 shifted down) needs its calculated-column formula restored.
$ws.Cells.Item(706, 7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# 5) Give the new row 522 the same look as the other "blank" data rows by
#    copying cell formatting from row 523 (which still has the original
#    formatting of the old row 522).
for ($c = 1; $c -le 11; $c++) {
    $ws.Cells.Item(523, $c).Copy()
    $ws.Cells.Item(522, $c).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# 5b) Column K on this new row holds a date, so pull the date-formatted
#     style from another "K" date cell (row 519) instead of the plain-text
#     style row 523 would have given it.
$ws.Cells.Item(519, 11).Copy()
$ws.Cells.Item(522, 11).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 6) Restore the calculated-column formula on the new row.
$ws.Cells.Item(522, 7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# 7) Populate the new leave-card entry: particulars + the date it was filed.
$ws.Cells.Item(522, 2).Value = "SOLO P(1-0-0)"
$ws.Cells.Item(522, 11).Value = 44991

# 8) Leave the selection where the author last clicked.
$ws.Range("B523").Select()
